$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "SCD0011"

# Update the TC_ID cell value
$ws.Range("B2").Value = "SCD0011-038"

# Update selection to B3
$ws.Range("B3").Select()

# Set column B width to match new content
$ws.Columns("B:B").ColumnWidth = 12.42578125
